# A new weekly data point was recorded for "Poroto verde" at "Feria
# Lagunitas de Puerto Montt". This pushes all existing rows from 70..82
# down to 71..83, and the freed-up row 70 is filled with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 70; Excel shifts rows 70-82 down to 71-83
# (matches xlShiftDown, the default for a whole-row insert).
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new weekly record.
$ws.Range("A70").Value = 4
$ws.Range("B70").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C70").Value = "Los Lagos"
$ws.Range("D70").Value = 44722
$ws.Range("E70").Value = 10
$ws.Range("F70").Value = 100112031
$ws.Range("G70").Value = "Poroto verde"
$ws.Range("H70").Value = "Magnum"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 40
$ws.Range("K70").Value = 27000
$ws.Range("L70").Value = 27000
$ws.Range("M70").Value = 27000
$ws.Range("N70").Value = "`$/malla 25 kilos"
$ws.Range("O70").Value = "Perú"
$ws.Range("P70").Value = 1080
$ws.Range("Q70").Value = 25
$ws.Range("R70").Value = "Hortaliza"
